$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.493404
$ws.Range("H2").Value = 46.48021199999999
$ws.Range("I2").Value = 0.05356331879335558
$ws.Range("J2").Value = 0.05356331879335557
$ws.Range("M2").Value = 0.8025546666666666
$ws.Range("N2").Value = 2.407664
$ws.Range("O2").Value = 0.1727979380778872
$ws.Range("P2").Value = 0.1727979380778872
$ws.Range("Q2").Value = 12.434303682752
$ws.Range("R2").Value = 111.908733144768
$ws.Range("S2").Value = 0.009255631044100389
$ws.Range("T2").Value = 0.009255631044100389
$ws.Range("G3").Value = 15.493404
$ws.Range("H3").Value = 46.48021199999999
$ws.Range("I3").Value = 0.05356331879335558
$ws.Range("J3").Value = 0.05356331879335557
$ws.Range("O3").Value = 0.4482183718598042
$ws.Range("P3").Value = 0.4482183718598043
$ws.Range("Q3").Value = 32.253182033812
$ws.Range("R3").Value = 290.2786383043079
$ws.Range("S3").Value = 0.02400806354096549
$ws.Range("T3").Value = 0.02400806354096549
$ws.Range("G4").Value = 15.493404
$ws.Range("H4").Value = 46.48021199999999
$ws.Range("I4").Value = 0.05356331879335558
$ws.Range("J4").Value = 0.05356331879335557
$ws.Range("M4").Value = 0.3607433333333334
$ws.Range("N4").Value = 1.08223
$ws.Range("O4").Value = 0.07767159891331675
$ws.Range("P4").Value = 0.07767159891331676
$ws.Range("Q4").Value = 5.58914220364
$ws.Range("R4").Value = 50.30227983275999
$ws.Range("S4").Value = 0.004160348613783636
$ws.Range("T4").Value = 0.004160348613783636
$ws.Range("G5").Value = 15.493404
$ws.Range("H5").Value = 46.48021199999999
$ws.Range("I5").Value = 0.05356331879335558
$ws.Range("J5").Value = 0.05356331879335557
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.399434666666667
$ws.Range("N5").Value = 4.198304
$ws.Range("O5").Value = 0.3013120911489918
$ws.Range("P5").Value = 0.3013120911489918
$ws.Range("Q5").Value = 21.682006662272
$ws.Range("R5").Value = 195.138059960448
$ws.Range("S5").Value = 0.01613927559450606
$ws.Range("T5").Value = 0.01613927559450606
$ws.Range("I6").Value = 0.4524333485785276
$ws.Range("J6").Value = 0.4524333485785275
$ws.Range("M6").Value = 0.8025546666666666
$ws.Range("N6").Value = 2.407664
$ws.Range("O6").Value = 0.1727979380778872
$ws.Range("P6").Value = 0.1727979380778872
$ws.Range("Q6").Value = 105.028847710752
$ws.Range("R6").Value = 945.259629396768
$ws.Range("S6").Value = 0.07817954975204355
$ws.Range("T6").Value = 0.07817954975204355
$ws.Range("I7").Value = 0.4524333485785276
$ws.Range("J7").Value = 0.4524333485785275
$ws.Range("O7").Value = 0.4482183718598042
$ws.Range("P7").Value = 0.4482183718598043
$ws.Range("S7").Value = 0.2027889388749469
$ws.Range("T7").Value = 0.2027889388749469
$ws.Range("I8").Value = 0.4524333485785276
$ws.Range("J8").Value = 0.4524333485785275
$ws.Range("M8").Value = 0.3607433333333334
$ws.Range("N8").Value = 1.08223
$ws.Range("O8").Value = 0.07767159891331675
$ws.Range("P8").Value = 0.07767159891331676
$ws.Range("Q8").Value = 47.20981410114
$ws.Range("R8").Value = 424.88832691026
$ws.Range("S8").Value = 0.03514122158580022
$ws.Range("T8").Value = 0.03514122158580022
$ws.Range("I9").Value = 0.4524333485785276
$ws.Range("J9").Value = 0.4524333485785275
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.399434666666667
$ws.Range("N9").Value = 4.198304
$ws.Range("O9").Value = 0.3013120911489918
$ws.Range("P9").Value = 0.3013120911489918
$ws.Range("Q9").Value = 183.141431470272
$ws.Range("R9").Value = 1648.272883232448
$ws.Range("S9").Value = 0.1363236383657369
$ws.Range("T9").Value = 0.1363236383657369
$ws.Range("G10").Value = 66.835223
$ws.Range("H10").Value = 200.505669
$ws.Range("I10").Value = 0.2310606730563543
$ws.Range("J10").Value = 0.2310606730563542
$ws.Range("M10").Value = 0.8025546666666666
$ws.Range("N10").Value = 2.407664
$ws.Range("O10").Value = 0.1727979380778872
$ws.Range("P10").Value = 0.1727979380778872
$ws.Range("Q10").Value = 53.63892011635733
$ws.Range("R10").Value = 482.750281047216
$ws.Range("S10").Value = 0.03992680787502684
$ws.Range("T10").Value = 0.03992680787502684
$ws.Range("G11").Value = 66.835223
$ws.Range("H11").Value = 200.505669
$ws.Range("I11").Value = 0.2310606730563543
$ws.Range("J11").Value = 0.2310606730563542
$ws.Range("O11").Value = 0.4482183718598042
$ws.Range("P11").Value = 0.4482183718598043
$ws.Range("Q11").Value = 139.1333120655356
$ws.Range("R11").Value = 1252.199808589821
$ws.Range("S11").Value = 0.1035656386781496
$ws.Range("T11").Value = 0.1035656386781496
$ws.Range("G12").Value = 66.835223
$ws.Range("H12").Value = 200.505669
$ws.Range("I12").Value = 0.2310606730563543
$ws.Range("J12").Value = 0.2310606730563542
$ws.Range("M12").Value = 0.3607433333333334
$ws.Range("N12").Value = 1.08223
$ws.Range("O12").Value = 0.07767159891331675
$ws.Range("P12").Value = 0.07767159891331676
$ws.Range("Q12").Value = 24.11036112909667
$ws.Range("R12").Value = 216.99325016187
$ws.Range("S12").Value = 0.01794685192227416
$ws.Range("T12").Value = 0.01794685192227416
$ws.Range("G13").Value = 66.835223
$ws.Range("H13").Value = 200.505669
$ws.Range("I13").Value = 0.2310606730563543
$ws.Range("J13").Value = 0.2310606730563542
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 1.399434666666667
$ws.Range("N13").Value = 4.198304
$ws.Range("O13").Value = 0.3013120911489918
$ws.Range("P13").Value = 0.3013120911489918
$ws.Range("Q13").Value = 93.53152802059735
$ws.Range("R13").Value = 841.783752185376
$ws.Range("S13").Value = 0.06962137458090362
$ws.Range("T13").Value = 0.06962137458090362
$ws.Range("G14").Value = 76.057215
$ws.Range("H14").Value = 228.171645
$ws.Range("I14").Value = 0.2629426595717627
$ws.Range("J14").Value = 0.2629426595717627
$ws.Range("M14").Value = 0.8025546666666666
$ws.Range("N14").Value = 2.407664
$ws.Range("O14").Value = 0.1727979380778872
$ws.Range("P14").Value = 0.1727979380778872
$ws.Range("Q14").Value = 61.04007283192
$ws.Range("R14").Value = 549.36065548728
$ws.Range("S14").Value = 0.04543594940671642
$ws.Range("T14").Value = 0.04543594940671643
$ws.Range("G15").Value = 76.057215
$ws.Range("H15").Value = 228.171645
$ws.Range("I15").Value = 0.2629426595717627
$ws.Range("J15").Value = 0.2629426595717627
$ws.Range("O15").Value = 0.4482183718598042
$ws.Range("P15").Value = 0.4482183718598043
$ws.Range("Q15").Value = 158.331067877645
$ws.Range("R15").Value = 1424.979610898805
$ws.Range("S15").Value = 0.1178557307657423
$ws.Range("T15").Value = 0.1178557307657423
$ws.Range("G16").Value = 76.057215
$ws.Range("H16").Value = 228.171645
$ws.Range("I16").Value = 0.2629426595717627
$ws.Range("J16").Value = 0.2629426595717627
$ws.Range("M16").Value = 0.3607433333333334
$ws.Range("N16").Value = 1.08223
$ws.Range("O16").Value = 0.07767159891331675
$ws.Range("P16").Value = 0.07767159891331676
$ws.Range("Q16").Value = 27.43713326315
$ws.Range("R16").Value = 246.93419936835
$ws.Range("S16").Value = 0.02042317679145874
$ws.Range("T16").Value = 0.02042317679145874
$ws.Range("G17").Value = 76.057215
$ws.Range("H17").Value = 228.171645
$ws.Range("I17").Value = 0.2629426595717627
$ws.Range("J17").Value = 0.2629426595717627
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 1.399434666666667
$ws.Range("N17").Value = 4.198304
$ws.Range("O17").Value = 0.3013120911489918
$ws.Range("P17").Value = 0.3013120911489918
$ws.Range("Q17").Value = 106.43710332112
$ws.Range("R17").Value = 957.9339298900801
$ws.Range("S17").Value = 0.07922780260784529
$ws.Range("T17").Value = 0.07922780260784529
